$d = $word.ActiveDocument
$d.Content.Find.Execute(" 0907107101", $true, $false, $false, $false, $false,
                         $true, 1, $false, " 0920038502", 2)
